$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 4: H4, I4, J4, K4, L4, M4, N4
$ws.Range("H4").Value = 1527.5555
$ws.Range("I4").Value = 229.6
$ws.Range("J4").Value = 3150
$ws.Range("K4").Value = 229.6
$ws.Range("L4").Value = 3150
$ws.Range("M4").Value = -115.6
$ws.Range("N4").Value = -3378
# Row 15: H15, I15, K15, M15
$ws.Range("H15").Value = 1163.78
$ws.Range("I15").Value = 1163.78
$ws.Range("K15").Value = 3491.34
$ws.Range("M15").Value = -3322.34
# Row 41: H41, I41, J41, K41, L41, M41, N41
$ws.Range("H41").Value = 883.875
$ws.Range("I41").Value = 967.75
$ws.Range("J41").Value = 800
$ws.Range("K41").Value = 967.75
$ws.Range("L41").Value = 800
$ws.Range("M41").Value = -527.75
$ws.Range("N41").Value = -1680
# Row 44: H44, I44, J44, K44, L44, M44, N44
$ws.Range("H44").Value = 8333
$ws.Range("I44").Value = 5000
$ws.Range("J44").Value = 9285.286
$ws.Range("K44").Value = 5000
$ws.Range("L44").Value = 9285.286
$ws.Range("M44").Value = -4538
$ws.Range("N44").Value = -10209.286
# Row 53: H53, I53, J53, K53, L53, M53, N53
$ws.Range("H53").Value = 179.9
$ws.Range("I53").Value = 151.25
$ws.Range("J53").Value = 199
$ws.Range("K53").Value = 151.25
$ws.Range("L53").Value = 199
$ws.Range("M53").Value = 485.75
$ws.Range("N53").Value = -1473
# Row 64: H64, J64, L64, N64
$ws.Range("H64").Value = 3550.8064
$ws.Range("J64").Value = 4350.75
$ws.Range("L64").Value = 4350.75
$ws.Range("N64").Value = -4846.75
# Row 67: H67, J67, L67, N67
$ws.Range("H67").Value = 3550.8064
$ws.Range("J67").Value = 4350.75
$ws.Range("L67").Value = 4350.75
$ws.Range("N67").Value = -6066.75
# Row 75: H75, J75, L75, N75
$ws.Range("H75").Value = 34352.5
$ws.Range("J75").Value = 34352.5
$ws.Range("L75").Value = 34352.5
$ws.Range("N75").Value = -36224.5
# Row 76: H76, I76, J76, K76, L76, M76, N76
$ws.Range("H76").Value = 3666.9048
$ws.Range("I76").Value = 3446.3333
$ws.Range("J76").Value = 4218.3335
$ws.Range("K76").Value = 3446.3333
$ws.Range("L76").Value = 4218.3335
$ws.Range("M76").Value = -3131.3333
$ws.Range("N76").Value = -4848.3335
# Row 78: H78, J78, L78, N78
$ws.Range("H78").Value = 34352.5
$ws.Range("J78").Value = 34352.5
$ws.Range("L78").Value = 103057.5
$ws.Range("N78").Value = -112417.5
# Row 79: H79, I79, J79, K79, L79, M79, N79
$ws.Range("H79").Value = 3666.9048
$ws.Range("I79").Value = 3446.3333
$ws.Range("J79").Value = 4218.3335
$ws.Range("K79").Value = 3446.3333
$ws.Range("L79").Value = 4218.3335
$ws.Range("M79").Value = -2354.3333
$ws.Range("N79").Value = -6402.3335
# Row 92: H92
$ws.Range("H92").Value = 622.1429000000001
# Row 116: H116, I116, J116, K116, L116, M116, N116
$ws.Range("H116").Value = 3555.5557
$ws.Range("I116").Value = 3166.6667
$ws.Range("J116").Value = 4333.3335
$ws.Range("K116").Value = 3166.6667
$ws.Range("L116").Value = 4333.3335
$ws.Range("M116").Value = 275.3332999999998
$ws.Range("N116").Value = -11217.3335
# Row 121: H121, I121, J121, K121, L121, M121, N121
$ws.Range("H121").Value = 1565.8334
$ws.Range("I121").Value = 1623.5
$ws.Range("J121").Value = 1450.5
$ws.Range("K121").Value = 4870.5
$ws.Range("L121").Value = 4351.5
$ws.Range("M121").Value = -3123.5
$ws.Range("N121").Value = -7845.5
# Row 137: H137, J137, L137, N137
$ws.Range("H137").Value = 3156.1924
$ws.Range("J137").Value = 8812
$ws.Range("L137").Value = 26436
$ws.Range("N137").Value = -31536

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61: H61, I61, K61, M61
$ws.Range("H61").Value = 3570.4443
$ws.Range("I61").Value = 2732.4243
$ws.Range("K61").Value = 2732.4243
$ws.Range("M61").Value = -2520.4243
# Row 74: H74, I74, J74, K74, L74, M74, N74
$ws.Range("H74").Value = 5903.7427
$ws.Range("I74").Value = 2416.3
$ws.Range("J74").Value = 26828.4
$ws.Range("K74").Value = 2416.3
$ws.Range("L74").Value = 26828.4
$ws.Range("M74").Value = -1542.3
$ws.Range("N74").Value = -28576.4
# Row 77: H77, I77, J77, K77, L77, M77, N77
$ws.Range("H77").Value = 5903.7427
$ws.Range("I77").Value = 2416.3
$ws.Range("J77").Value = 26828.4
$ws.Range("K77").Value = 12081.5
$ws.Range("L77").Value = 134142
$ws.Range("M77").Value = -7713.5
$ws.Range("N77").Value = -142878
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 7732.2104
$ws.Range("I132").Value = 6713.0527
$ws.Range("J132").Value = 8751.368
$ws.Range("K132").Value = 20139.1581
$ws.Range("L132").Value = 26254.104
$ws.Range("M132").Value = -17609.1581
$ws.Range("N132").Value = -31314.104
# Row 133: H133, J133, L133, N133
$ws.Range("H133").Value = 20328.766
$ws.Range("J133").Value = 20328.766
$ws.Range("L133").Value = 20328.766
$ws.Range("N133").Value = -25388.766
# Row 136: H136, I136, K136, M136
$ws.Range("H136").Value = 3570.4443
$ws.Range("I136").Value = 2732.4243
$ws.Range("K136").Value = 8197.2729
$ws.Range("M136").Value = -5647.2729

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134: H134, I134, J134, K134, L134, M134, N134
$ws.Range("H134").Value = 8314.947
$ws.Range("I134").Value = 9463.4375
$ws.Range("J134").Value = 2189.6667
$ws.Range("K134").Value = 28390.3125
$ws.Range("L134").Value = 6569.000100000001
$ws.Range("M134").Value = -25855.3125
$ws.Range("N134").Value = -11639.0001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22: H22, I22, K22, M22
$ws.Range("H22").Value = 220.5
$ws.Range("I22").Value = 235.66667
$ws.Range("K22").Value = 235.66667
$ws.Range("M22").Value = 114.33333
# Row 31: H31, I31, J31, K31, L31, M31, N31
$ws.Range("H31").Value = 4216.143
$ws.Range("I31").Value = 4008.5532
$ws.Range("J31").Value = 5300.222
$ws.Range("K31").Value = 4008.5532
$ws.Range("L31").Value = 5300.222
$ws.Range("M31").Value = -3713.5532
$ws.Range("N31").Value = -5890.222
# Row 34: H34, I34, J34, K34, L34, M34, N34
$ws.Range("H34").Value = 4216.143
$ws.Range("I34").Value = 4008.5532
$ws.Range("J34").Value = 5300.222
$ws.Range("K34").Value = 4008.5532
$ws.Range("L34").Value = 5300.222
$ws.Range("M34").Value = -3806.5532
$ws.Range("N34").Value = -5704.222
# Row 58: H58, I58, J58, K58, L58, M58, N58
$ws.Range("H58").Value = 2528520.2
$ws.Range("I58").Value = 5349776.5
$ws.Range("J58").Value = 4238.316
$ws.Range("K58").Value = 5349776.5
$ws.Range("L58").Value = 4238.316
$ws.Range("M58").Value = -5349573.5
$ws.Range("N58").Value = -4644.316
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 1831.6487
$ws.Range("I132").Value = 1501.7307
$ws.Range("J132").Value = 2611.4546
$ws.Range("K132").Value = 4505.1921
$ws.Range("L132").Value = 7834.3638
$ws.Range("M132").Value = -1975.1921
$ws.Range("N132").Value = -12894.3638
# Row 134: H134, I134, K134, M134
$ws.Range("H134").Value = 2901.0454
$ws.Range("I134").Value = 2247.4167
$ws.Range("K134").Value = 6742.250100000001
$ws.Range("M134").Value = -4207.250100000001
# Row 136: H136, I136, J136, K136, L136, M136, N136
$ws.Range("H136").Value = 2528520.2
$ws.Range("I136").Value = 5349776.5
$ws.Range("J136").Value = 4238.316
$ws.Range("K136").Value = 16049329.5
$ws.Range("L136").Value = 12714.948
$ws.Range("M136").Value = -16046779.5
$ws.Range("N136").Value = -17814.948

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 60: H60, I60, J60, K60, L60, M60, N60
$ws.Range("H60").Value = 815
$ws.Range("I60").Value = 650
$ws.Range("J60").Value = 980
$ws.Range("K60").Value = 1950
$ws.Range("L60").Value = 2940
$ws.Range("M60").Value = -1699
$ws.Range("N60").Value = -3442

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97: H97, I97, J97, K97, L97, M97, N97
$ws.Range("H97").Value = 1381
$ws.Range("I97").Value = 1153.3636
$ws.Range("J97").Value = 1659.2222
$ws.Range("K97").Value = 1153.3636
$ws.Range("L97").Value = 1659.2222
$ws.Range("M97").Value = -657.3635999999999
$ws.Range("N97").Value = -2651.2222
# Row 103: H103, J103, L103, N103
$ws.Range("H103").Value = 45000
$ws.Range("J103").Value = 45000
$ws.Range("L103").Value = 45000
$ws.Range("N103").Value = -47344
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 2637.75
$ws.Range("I132").Value = 2131
$ws.Range("J132").Value = 2999.7144
$ws.Range("K132").Value = 6393
$ws.Range("L132").Value = 8999.143199999999
$ws.Range("M132").Value = -3863
$ws.Range("N132").Value = -14059.1432

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22: H22, I22, J22, K22, L22, M22, N22
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 633.3333
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 633.3333
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -338.3333
$ws.Range("N22").Value = -1090
# Row 27: H27, I27, J27, K27, L27, M27, N27
$ws.Range("H27").Value = 600
$ws.Range("I27").Value = 633.3333
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 633.3333
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -526.3333
$ws.Range("N27").Value = -714

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81: H81, I81, J81, K81, L81, M81, N81
$ws.Range("H81").Value = 1768.1818
$ws.Range("I81").Value = 1010
$ws.Range("J81").Value = 2400
$ws.Range("K81").Value = 2020
$ws.Range("L81").Value = 4800
$ws.Range("M81").Value = -959
$ws.Range("N81").Value = -6922
# Row 84: H84, I84, J84, K84, L84, M84, N84
$ws.Range("H84").Value = 1768.1818
$ws.Range("I84").Value = 1010
$ws.Range("J84").Value = 2400
$ws.Range("K84").Value = 10100
$ws.Range("L84").Value = 24000
$ws.Range("M84").Value = -4796
$ws.Range("N84").Value = -34608
# Row 101: H101, J101, L101, N101
$ws.Range("H101").Value = 54000
$ws.Range("J101").Value = 54000
$ws.Range("L101").Value = 54000
$ws.Range("N101").Value = -60490
# Row 126: H126, I126, J126, K126, L126, M126, N126
$ws.Range("H126").Value = 1238.1428
$ws.Range("I126").Value = 1225.6875
$ws.Range("J126").Value = 1278
$ws.Range("K126").Value = 3677.0625
$ws.Range("L126").Value = 3834
$ws.Range("M126").Value = -1207.0625
$ws.Range("N126").Value = -8774
# Row 132: H132, I132, J132, K132, L132, M132, N132
$ws.Range("H132").Value = 1523.325
$ws.Range("I132").Value = 787.16
$ws.Range("J132").Value = 2750.2666
$ws.Range("K132").Value = 2361.48
$ws.Range("L132").Value = 8250.799800000001
$ws.Range("M132").Value = 168.52
$ws.Range("N132").Value = -13310.7998
